$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.474.77"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "2.362.40"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'521.76"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'136.33"
$ws.Range("E6").Value = "  +1.64%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").Value = "'5.44"
$ws.Range("E10").Value = "  +5.59%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "'0.343"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "'24.31"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").Value = "2.785.57"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").Value = "57.527.19"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "2.378.33"
$ws.Range("E17").Value = "  +2.72%  "
$ws.Range("D18").Value = "'10.64"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").Value = "'330.84"
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").Value = "'6.72"
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'61.35"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'8.79"
$ws.Range("E24").Value = "  +14.34%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "'0.166"
$ws.Range("E25").Value = "  +4.55%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "'1.35"
$ws.Range("E27").Value = "  +11.81%  "
$ws.Range("D28").Value = "0.0₃0747"
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("D29").Value = "'169.89"
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").Value = "'6.31"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'18.60"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +3.65%  "
$ws.Range("D35").Value = "'0.994"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").Value = "'0.925"
$ws.Range("E36").Value = "  -2.61%  "
$ws.Range("D37").Value = "'4.06"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").Value = "'1.64"
$ws.Range("E38").Value = "  +7.64%  "
$ws.Range("D39").Value = "'38.67"
$ws.Range("E39").Value = "  +3.06%  "
$ws.Range("D40").Value = "'151.43"
$ws.Range("E40").Value = "  +7.87%  "
$ws.Range("D41").Value = "'0.387"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("D43").Value = "'5.36"
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("D44").Value = "'283.75"
$ws.Range("E44").Value = "  +2.85%  "
$ws.Range("D45").Value = "'0.0943"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'0.567"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'18.36"
$ws.Range("E48").Value = "  +6.51%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0221"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("D50").Value = "'17.81"
$ws.Range("E50").Value = "  +5.17%  "
$ws.Range("D51").Value = "'1.66"
$ws.Range("E51").Value = "  +0.39%  "
